$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Foundation_PConc",
    "Foundation_CBlock",
    "GarageCars",
    "Fireplaces",
    "TotalBsmtSF",
    "GarageFinish_Unf",
    "Neighborhood_NridgHt",
    "MSZoning_RM",
    "BsmtExposure_Gd",
    "HeatingQC",
    "Neighborhood_NoRidge",
    "SaleType_New",
    "GarageFinish_No Garage",
    "BsmtFinType1_GLQ",
    "GarageType_Attchd",
    "GarageType_BuiltIn",
    "SaleCondition_Partial",
    "RoofStyle_Gable",
    "BuildingAge",
    "Exterior2nd_VinylSd",
    "Neighborhood_NAmes",
    "LotArea",
    "FullBath",
    "CentralAir",
    "MSZoning_RL",
    "MasVnrArea",
    "BsmtExposure_No",
    "BsmtFinSF1",
    "SaleType_WD",
    "Neighborhood_MeadowV",
    "GarageType_Detchd",
    "OverallCond",
    "HouseStyle_2Story",
    "Exterior1st_BrkFace",
    "RoofStyle_Hip",
    "LotFrontage",
    "BsmtExposure_No Bsmnt",
    "Neighborhood_StoneBr",
    "Neighborhood_Crawfor",
    "Exterior1st_CemntBd",
    "GarageQual",
    "GarageType_No Garage",
    "2ndFlrSF",
    "BsmtCond"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
